# "i have update the logout method"
#
# The Test Suite sheet drives which regression suites execute. The
# RegressionSuite row (row 2) - which exercises the logout flow - was
# switched on: its Runmode cell (C2) flips from "n" to "y".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2").Value = "y"

# Normalize formatting on the APISuite row (row 5): drop the redundant
# "apply fill" override left over from earlier edits so it matches the
# plain bordered style already used by the other data rows.
$xlNone = -4142
$ws.Range("A5").Interior.Pattern = $xlNone
$ws.Range("C5").Interior.Pattern = $xlNone
